$wb = $excel.ActiveWorkbook

# --- Sheet "2_": was the active tab with B4 selected; now deactivated and
#     reselected as the full A1:C5 range (its content is unchanged). ---
$ws2 = $wb.Worksheets.Item("2_")
$ws2.Activate()
$ws2.Range("A1:C5").Select()

# --- Sheet "3_": becomes the new active tab, gains the Zeroth-law True/False
#     question content (mirroring the "True_False" question bank), and the
#     leftover placeholder answers in B3:B5 are cleared out. ---
$ws3 = $wb.Worksheets.Item("3_")
$ws3.Activate()

$ws3.Range("A1").Value = "Which of the following must be true, according to the Zeroth law?  Mark each with a 'T' or 'F'."
$ws3.Range("A2").Value = "Two objects that have the same temperature must be at equilibrium"
$ws3.Range("B2").Value = "T"
$ws3.Range("C2").Value = "This is true: the Zeroth law stipulates that equal temperature is in indication of equilibrium"

$ws3.Range("B3").Value = ""
$ws3.Range("B4").Value = ""
$ws3.Range("B5").Value = ""

$ws3.Rows.Item(1).RowHeight = 45
$ws3.Rows.Item(2).RowHeight = 45

$ws3.Range("C10").Select()
